$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.050.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.638.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.733.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.070.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "191.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.132"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0487"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.77%  "
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.131.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.785"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "55.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0528"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E51").Value = "  -0.40%  "